$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "왜 코딩하는 분들은 각 데이터마다 적합한 모델이 있다는걸 이해를 못하는걸까요?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/why-coders-are-stupid/#utm_source=rss&utm_medium=rss&utm_campaign=why-coders-are-stupid"

$ws.Range("D29").Value = "[만화] 인턴일기 58~65"

$ws.Range("D36").Value = "Handling imbalanced datasets"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/343"

$ws.Range("D50").Value = "vanishing gradient problem"
$ws.Range("E50").Value = "http://incredible.egloos.com/7530413"
